$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.323.70"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.667.75"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'220.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'0.5308"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.2646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "'0.06360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "'20.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "'0.07836"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'4.529"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "1.670.11"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "1.895.87"
$ws.Range("D15").Value = "'0.5610"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "0.0₅8136"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "'65.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "26.313.20"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'4.718"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "'198.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "'6.057"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'146.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").Value = "'7.231"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'16.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "'1.518"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.30%  "
$ws.Range("D30").Value = "'0.05889"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "'3.538"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "'3.321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'2.832"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.9614"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").Value = "'2.432"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.5804"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").Value = "'5.967"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").Value = "1.074.88"
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("D42").Value = "'0.8564"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'58.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4415"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.080"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈103"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").Value = "'0.05152"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.16%  "
